# Update "C" column values in result_data_RandomForest.xlsx (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value  = -11.3796
$ws.Range("C21").Value = -12.62680000000001
$ws.Range("C23").Value = -12.3379
$ws.Range("C25").Value = -13.3261
$ws.Range("C53").Value = -10.50860000000001
$ws.Range("C57").Value = -13.90739999999999
$ws.Range("C59").Value = -12.9065
$ws.Range("C69").Value = -11.90390000000001
$ws.Range("C79").Value = -10.31870000000001
$ws.Range("C83").Value = -13.65569999999999
$ws.Range("C93").Value = -11.0334
